# 2021-06 QLD Outbreak Paths - "Add files via upload"
#
# A new outbreak-chain edge was added to Sheet1 (Q23 -> Q24 w, dated
# 2021-07-06, Queensland cluster, Alpha variant, Isolated), which extends
# Table1 by one row. The "Date Colours" helper sheet got a matching new
# date-colour entry (row 18, 2021-07-06) and its colour gradient (column B,
# and the helper spill row 2 F:V) shifted down to make room for it.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Sheet1"            -> Table1
$ws2 = $wb.Worksheets.Item(2)   # "Date Colours"       -> Date_Colours

# ---------------------------------------------------------------------
# Sheet1: append new outbreak-path row 29
#   Date | Source | Target | Cluster | Sub-Cluster | Link Label | Variant | Isolated
# ---------------------------------------------------------------------
$ws1.Range("A29").Value = 44383
$ws1.Range("A29").NumberFormat = $ws1.Range("A28").NumberFormat
$ws1.Range("B29").Value = "Q23"
$ws1.Range("C29").Value = "Q24 w"
$ws1.Range("D29").Value = "Queensland"
$ws1.Range("G29").Value = "Alpha (B.1.1.7)"
$ws1.Range("H29").Value = "Isolated"

# Grow Table1 to include the new row
$table1 = $ws1.ListObjects.Item(1)
$table1.Resize($ws1.Range("A1:H29"))

# ---------------------------------------------------------------------
# Date Colours sheet: shift the purple colour gradient down one row
# (column B, rows 2-18) and refresh the helper gradient row (F2:V2)
# ---------------------------------------------------------------------
$ws2.Range("B2").Value  = "#faf5fa"
$ws2.Range("B3").Value  = "#f4ebf4"
$ws2.Range("B4").Value  = "#efe1ef"
$ws2.Range("B5").Value  = "#e9d6ea"
$ws2.Range("B6").Value  = "#e4cce5"
$ws2.Range("B7").Value  = "#dec3df"
$ws2.Range("B8").Value  = "#d8b9da"
$ws2.Range("B9").Value  = "#d3afd5"
$ws2.Range("B10").Value = "#cda5d0"
$ws2.Range("B11").Value = "#c79bca"
$ws2.Range("B12").Value = "#c291c5"
$ws2.Range("B13").Value = "#bc88c0"
$ws2.Range("B14").Value = "#b67ebb"
$ws2.Range("B15").Value = "#b074b6"
$ws2.Range("B16").Value = "#aa6bb0"
$ws2.Range("B17").Value = "#a461ab"
$ws2.Range("B18").Value = "#9e57a6"

$ws2.Range("F2").Value  = "#faf5fa"
$ws2.Range("G2").Value  = "#f4ebf4"
$ws2.Range("H2").Value  = "#efe1ef"
$ws2.Range("I2").Value  = "#e9d6ea"
$ws2.Range("J2").Value  = "#e4cce5"
$ws2.Range("K2").Value  = "#dec3df"
$ws2.Range("L2").Value  = "#d8b9da"
$ws2.Range("M2").Value  = "#d3afd5"
$ws2.Range("N2").Value  = "#cda5d0"
$ws2.Range("O2").Value  = "#c79bca"
$ws2.Range("P2").Value  = "#c291c5"
$ws2.Range("Q2").Value  = "#bc88c0"
$ws2.Range("R2").Value  = "#b67ebb"
$ws2.Range("S2").Value  = "#b074b6"
$ws2.Range("T2").Value  = "#aa6bb0"
$ws2.Range("U2").Value  = "#a461ab"
$ws2.Range("V2").Value  = "#9e57a6"

# ---------------------------------------------------------------------
# Restore on-screen selections left by the editor (Sheet2 selected first
# so Sheet1 ends up the active/selected tab, matching the saved file).
# ---------------------------------------------------------------------
$ws2.Range("B2:B18").Select()
$ws1.Range("F29").Select()
